$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update source-of-data text: "государственной" -> "официальной"
$ws.Range("B16").Value = "Источниками информации для расчета показателя является ежегодные данные официальной статистической отчетности, представляемой дошкольными организациями республики по форме «Отчет дошкольной организации», а также общеобразовательными организациями республики по форме ОШ-1 «Отчет дневной общеобразовательной школы на начало учебного года»."

# Update references text: law date/number changed
$ws.Range("B26").Value = "https://sustainabledevelopment-kyrgyzstan.github.io`nwww.uis.unesco.org `nЗакон Кыргызской Республики от августа 2023 года №179 «Об образовании»"

# Update view: zoom and active selection
$ws.Application.ActiveWindow.Zoom = 82
$ws.Range("B2").Select()
